$wb = $excel.ActiveWorkbook

# --- Rename sheets (sheetId order: GNG_TO, NB_TO, RS_TO, TOL_TO, vSAT_TO) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509961993064287"
$wb.Worksheets.Item(2).Name = "NB_TO-16509962010264375"
$wb.Worksheets.Item(3).Name = "RS_TO-16509962010264375"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509962010744338"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509962011384366"

# --- Sheet 1 (GNG_TO) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961992664.csv"
$ws1.Range("B3").Value = "GNG_stims-1650996199290432.csv"
$ws1.Range("B4").Value = "go_stims-1650996199290432.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961993064287.csv"

# --- Sheet 2 (NB_TO) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16509962003383956.csv"
$ws2.Range("B3").Value = "OB-16509962004183986.csv"
$ws2.Range("B4").Value = "ZB-match_1-16509961997303967.csv"
$ws2.Range("B5").Value = "TB-1650996200890432.csv"
$ws2.Range("B6").Value = "ZB-match_5-1650996200050438.csv"
$ws2.Range("B7").Value = "TB-16509962010023983.csv"
$ws2.Range("B8").Value = "ZB-match_3-16509961997704012.csv"
$ws2.Range("B9").Value = "TB-1650996200626404.csv"
$ws2.Range("B10").Value = "OB-16509962004504342.csv"

# --- Sheet 3 (RS_TO) --- only the sheet name changed, no cell content changes

# --- Sheet 4 (TOL_TO) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509962010424006.csv"
$ws4.Range("B3").Value = "ZM_stims-16509962010264375.csv"
$ws4.Range("B4").Value = "MM_stims-1650996201058439.csv"
$ws4.Range("B5").Value = "ZM_stims-16509962010424006.csv"
$ws4.Range("B6").Value = "MM_stims-16509962010744338.csv"
$ws4.Range("B7").Value = "ZM_stims-1650996201058439.csv"

# --- Sheet 5 (vSAT_TO) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509962011064448.csv"
$ws5.Range("B3").Value = "SAT_stims-1650996201090411.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509962011224008.csv"
$ws5.Range("B5").Value = "SAT_stims-16509962010744338.csv"
